$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.809.64"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.807.09"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.67"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.20"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.453"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.08"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "4.446.77"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "3.820.84"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.836.74"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.46"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.19"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000148"
$ws.Range("E23").Value = "  -3.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.31"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "3.956.94"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.48"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.82"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.27"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.301"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.90"
$ws.Range("E46").Value = "  +8.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.73"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +11.45%  "
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "391.19"
$ws.Range("E51").Value = "  +0.34%  "
